# Appends the git-workflow "cheat sheet" lines to the end of the document,
# matching the author's "j ai ajouté encore" commit.
#
# Strategy: first create all the new (still-empty/plain) paragraphs via
# InsertParagraphAfter, *then* go back and fill in text/formatting on each
# one individually. Doing the formatting only after every paragraph already
# exists avoids "leaking" a paragraph's character formatting forward onto
# the next paragraph that InsertParagraphAfter() would otherwise clone it
# into.

$d = $word.ActiveDocument

# The document currently ends with an empty paragraph. Remember its 1-based
# index, then append 10 more (still empty/plain) paragraphs after it, for
# 11 new paragraphs total. Every paragraph is created *before* any text or
# character formatting is applied below -- InsertParagraphAfter() clones
# the current paragraph-mark formatting onto the new paragraph, so doing
# all the structural inserts first (while formatting is still plain)
# avoids leaking formatting from one new line into the next.
$firstNewIndex = $d.Paragraphs.Count
for ($i = 0; $i -lt 10; $i++) {
    $last = $d.Paragraphs.Last
    $last.Range.InsertParagraphAfter()
}

# --- paragraph firstNewIndex: existing last/empty paragraph gets its text -
$pAdd = $d.Paragraphs.Item($firstNewIndex)
$pAdd.Range.InsertAfter('git add ligne_commandes                 ou           git add .')

# --- next two paragraphs: left blank ---------------------------------------

# --- "git commit -m ..." : monospace / black text / white highlight --------
$pCommit = $d.Paragraphs.Item($firstNewIndex + 3)
$pCommit.Range.InsertAfter('git commit -m "Ajouté ma checklist-vacances.md"')
$rCommit = $pCommit.Range
$rCommit.Font.Name = "monospace;monospace"
$rCommit.Font.Bold = $false
$rCommit.Font.Italic = $false
$rCommit.Font.AllCaps = $false
$rCommit.Font.SmallCaps = $false
$rCommit.Font.Color = 0
$rCommit.Font.Spacing = 0
$rCommit.Font.Size = 10.5
$rCommit.HighlightColorIndex = 8

# --- "git log ..." -----------------------------------------------------------
$pLog = $d.Paragraphs.Item($firstNewIndex + 4)
$pLog.Range.InsertAfter('git log                : affiche tous les commit et leur SHA')

# --- next paragraph: left blank ---------------------------------------------

# --- "git commit -a -m ..." --------------------------------------------------
$pCommitA = $d.Paragraphs.Item($firstNewIndex + 6)
$pCommitA.Range.InsertAfter('git commit -a -m "j ai ajouté un commit »')

# --- next two paragraphs: left blank -----------------------------------------

# --- "git checkout SHADuCommit" ----------------------------------------------
$pCheckout = $d.Paragraphs.Item($firstNewIndex + 9)
$pCheckout.Range.InsertAfter('git checkout SHADuCommit')

# --- trailing empty paragraph: monospace / white text / black highlight -----
$pTail = $d.Paragraphs.Item($firstNewIndex + 10)
$rTail = $pTail.Range
$rTail.Font.Name = "monospace;monospace"
$rTail.Font.Bold = $false
$rTail.Font.Italic = $false
$rTail.Font.AllCaps = $false
$rTail.Font.SmallCaps = $false
$rTail.Font.Color = 16777215
$rTail.Font.Spacing = 0
$rTail.Font.Size = 10.5
$rTail.HighlightColorIndex = 1
